$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab from "Sheet1" to "raw data"
$ws.Name = "raw data"

# Fix the column D header text (shared string "Answer_relevance" -> "Answer relevance")
$ws.Range("D1").Value = "Answer relevance"

# Row height adjustments (wrapped-text rows re-measured to slightly shorter heights)
$ws.Rows.Item(4).RowHeight = 240
$ws.Rows.Item(6).RowHeight = 304
$ws.Rows.Item(7).RowHeight = 320
$ws.Rows.Item(8).RowHeight = 350
$ws.Rows.Item(10).RowHeight = 208
$ws.Rows.Item(14).RowHeight = 256
$ws.Rows.Item(15).RowHeight = 96
$ws.Rows.Item(18).RowHeight = 240
$ws.Rows.Item(19).RowHeight = 256
$ws.Rows.Item(21).RowHeight = 335
$ws.Rows.Item(22).RowHeight = 335
$ws.Rows.Item(25).RowHeight = 240
$ws.Rows.Item(26).RowHeight = 395
$ws.Rows.Item(28).RowHeight = 256
$ws.Rows.Item(30).RowHeight = 208
$ws.Rows.Item(35).RowHeight = 320
$ws.Rows.Item(36).RowHeight = 320
$ws.Rows.Item(39).RowHeight = 256
$ws.Rows.Item(42).RowHeight = 256
$ws.Rows.Item(46).RowHeight = 208
$ws.Rows.Item(49).RowHeight = 335
$ws.Rows.Item(50).RowHeight = 304
$ws.Rows.Item(55).RowHeight = 335
$ws.Rows.Item(56).RowHeight = 256
$ws.Rows.Item(61).RowHeight = 256
$ws.Rows.Item(63).RowHeight = 320
$ws.Rows.Item(64).RowHeight = 335
$ws.Rows.Item(67).RowHeight = 256
$ws.Rows.Item(70).RowHeight = 256

# Restore view to the top of the sheet and select D2 (previously scrolled to row 70, cell D71 selected)
$ws.Range("D2").Select()
